$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column E (Mean Sign Production) before the current
# "% Developmental Delays" column, shifting it to F.
$ws.Range("E1").EntireColumn.Insert()

# Update header row
$ws.Range("D1").Value = "Mean Spoken Word Production (SD)"
$ws.Range("E1").Value = "Mean Sign Production (SD)"
$ws.Range("F1").Value = "% Developmental Delays"

# Row 2 - WG group
$ws.Range("A2").Value = "WG (n = 75)"
$ws.Range("B2").Value = "20.1 (8.9) months"
$ws.Range("C2").Value = "106 (99) words"
$ws.Range("D2").Value = "36 (59) words"
$ws.Range("E2").Value = "0 (2) words"
$ws.Range("F2").Value = "'18.7%"
$ws.Range("F2").Style = "Normal"

# Row 3 - WS group
$ws.Range("A3").Value = "WS (n = 24)"
$ws.Range("B3").Value = "25.8 (7.7) months"
$ws.Range("C3").Value = "NA"
$ws.Range("D3").Value = "138 (185) words"
$ws.Range("E3").Value = "0 (0) words"
$ws.Range("F3").Value = "'4.2%"
$ws.Range("F3").Style = "Normal"
